$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width: source diff goes from 16.42578125 to 15.42578125 (delta -1 character unit).
# This COM runtime snaps ColumnWidth to the nearest 1/6 (pixel-based) increment, so an input
# of 14.67 lands on the closest reachable stored width (15.5) to the target 15.42578125.
$ws.Columns.Item(1).ColumnWidth = 14.67

# New values for A1:A33 (column vector of doubles)
$newValues = @(
    0.18419135961709543,
    -0.0059999999397248871,
    -0.003999999945806465,
    -0.0079999999006901135,
    -0.0029999999425811552,
    -0.0019999999340409857,
    -0.0099999998552031677,
    -0.0099999998529951561,
    0.0043612832179240613,
    -0.031661300470071652,
    -0.0029999999172858338,
    0.043250399183788257,
    -0.0034999999072695687,
    -0.0079999998618998092,
    -0.00099999992898602841,
    -0.0019999999183037964,
    -0.0019999999169417748,
    -0.0039999998972435336,
    -0.0039999999589985791,
    -0.0039999999554982679,
    -0.0039999999550026644,
    -0.0039999999548285814,
    -0.0049999999325809341,
    -0.019999999780016431,
    -0.019999999777065902,
    -0.0024999999271191342,
    -0.0024999999239989634,
    -0.0019999999152773285,
    -0.0069999998570891009,
    -0.059999999337757348,
    0.056032308425193023,
    -0.04486670584768504,
    -0.0039999998799356007
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $newValues[$i]
}
